$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

$c = $t.Cell(1, 1)
$c.Range.Text = "89 x 18" + $vt + "  1    8" + $vt + "  ----" + $vt + "8|    |" + $vt + "9|    |"
$c = $t.Cell(1, 2)
$c.Range.Text = "56 x 62" + $vt + "  6    2" + $vt + "  ----" + $vt + "5|    |" + $vt + "6|    |"
$c = $t.Cell(1, 3)
$c.Range.Text = "65 x 78" + $vt + "  7    8" + $vt + "  ----" + $vt + "6|    |" + $vt + "5|    |"
$c = $t.Cell(2, 1)
$c.Range.Text = "49 x 25" + $vt + "  2    5" + $vt + "  ----" + $vt + "4|    |" + $vt + "9|    |"
$c = $t.Cell(2, 2)
$c.Range.Text = "41 x 93" + $vt + "  9    3" + $vt + "  ----" + $vt + "4|    |" + $vt + "1|    |"
$c = $t.Cell(2, 3)
$c.Range.Text = "34 x 82" + $vt + "  8    2" + $vt + "  ----" + $vt + "3|    |" + $vt + "4|    |"
$c = $t.Cell(3, 1)
$c.Range.Text = "86 x 95" + $vt + "  9    5" + $vt + "  ----" + $vt + "8|    |" + $vt + "6|    |"
$c = $t.Cell(3, 2)
$c.Range.Text = "12 x 16" + $vt + "  1    6" + $vt + "  ----" + $vt + "1|    |" + $vt + "2|    |"
$c = $t.Cell(3, 3)
$c.Range.Text = "82 x 39" + $vt + "  3    9" + $vt + "  ----" + $vt + "8|    |" + $vt + "2|    |"
$c = $t.Cell(4, 1)
$c.Range.Text = "21 x 49" + $vt + "  4    9" + $vt + "  ----" + $vt + "2|    |" + $vt + "1|    |"
$c = $t.Cell(4, 2)
$c.Range.Text = "58 x 36" + $vt + "  3    6" + $vt + "  ----" + $vt + "5|    |" + $vt + "8|    |"
$c = $t.Cell(4, 3)
$c.Range.Text = "36 x 21" + $vt + "  2    1" + $vt + "  ----" + $vt + "3|    |" + $vt + "6|    |"
$c = $t.Cell(5, 1)
$c.Range.Text = "80 x 15" + $vt + "  1    5" + $vt + "  ----" + $vt + "8|    |" + $vt + "0|    |"
$c = $t.Cell(5, 2)
$c.Range.Text = "63 x 91" + $vt + "  9    1" + $vt + "  ----" + $vt + "6|    |" + $vt + "3|    |"
$c = $t.Cell(5, 3)
$c.Range.Text = "17 x 46" + $vt + "  4    6" + $vt + "  ----" + $vt + "1|    |" + $vt + "7|    |"
